# Rename "Attribut..." shared strings to "Attribute..." for consistency,
# and reorder the IfcSpace PredefinedType / IsInteriorOrExteriorSpace /
# IsExternal attribute rows (with the Pset_SpaceCommon marker moving from
# row 10 to row 9) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AttributName -> AttributeName, AttributDescriptionFR -> AttributeDescriptionFR
$ws.Range("D1:E1").Value = "AttributeName"
$ws.Range("F1").Value = "AttributeDescriptionFR"

# IfcSpace attribute rows 8-10 get reordered:
#  row 8: IsExternal               -> PredefinedType
#  row 9: PredefinedType           -> IsInteriorOrExteriorSpace (gains Pset_SpaceCommon in C)
#  row 10: IsInteriorOrExteriorSpace -> IsExternal (loses Pset_SpaceCommon in C)
$ws.Range("D8:E8").Value = "PredefinedType"
$ws.Range("D9:E9").Value = "IsInteriorOrExteriorSpace"
$ws.Range("D10:E10").Value = "IsExternal"

$ws.Range("C9").Value = "Pset_SpaceCommon"
$ws.Range("C10").Clear()

# Column width adjustments: column C narrows to match A/B/D, H widens to
# match G/I, and J/K narrow down to the same width as the L:V block.
$ws.Columns.Item(3).ColumnWidth = 19.75
$ws.Columns.Item(8).ColumnWidth = 19.75
$ws.Columns.Item(10).ColumnWidth = 7.75
$ws.Columns.Item(11).ColumnWidth = 7.75
